{"js": "// Modificado descripcion RF registro usuario\n//\n// Appends a new sentence about BCrypt password encryption to the end of\n// the \"Registro de usuario\" (user registration) functional-requirement\n// paragraph, matching the run formatting (sz=24 / szCs=24, i.e. 12pt)\n// used by the rest of the paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Registro de usuario\" paragraph that talks about the client\n// registering with name/email/password/birth date/municipality/address\n// (there are two similar paragraphs in the doc; this one is identified by\n// the trailing \"municipio y direcci\u00f3n\" wording that is unique to it).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Registro de usuario\") !== -1 &&\n      text.indexOf(\"municipio y direcci\u00f3n\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('No se encontr\u00f3 el p\u00e1rrafo de \"Registro de usuario\".');\n}\n\nconst sentence =\n  \" As\u00ed mismo, la contrase\u00f1a del usuario ser\u00e1 encriptada utilizando la \" +\n  \"librer\u00eda BCrypt que proporciona Java, utilizando una \\u2018sal\\u2019 \" +\n  \"para reforzar a\u00fan m\u00e1s la contrase\u00f1a a la hora de almacenarla en la \" +\n  \"base de datos.\";\n\n// Use insertOoxml so the new run carries the exact same run properties\n// (sz/szCs = 24, i.e. 12pt) as the rest of the paragraph's body text.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">' + sentence + '</w:t></w:r></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Modificado descripcion RF registro usuario\n#\n# Appends a new sentence about BCrypt password encryption to the end of\n# the \"Registro de usuario\" (user registration) functional-requirement\n# paragraph, matching the run formatting (sz=24 / szCs=24, i.e. 12pt)\n# used by the rest of the paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the end of the \"... municipio y direcci\u00f3n.\" sentence - this text\n# is unique to the specific \"Registro de usuario\" paragraph that needs the\n# new sentence appended (there is a similar, but textually different,\n# \"Registro de usuario\" paragraph elsewhere in the document that must stay\n# untouched).\n$anchor = $d.Content\n$found = $anchor.Find.Execute(\"municipio y direcci\u00f3n.\")\nif (-not $found) {\n    throw 'No se encontr\u00f3 el texto ancla \"municipio y direcci\u00f3n.\" en el documento.'\n}\n\n# Collapse the found range to its end (right after the period) and\n# remember that insertion point.\n$anchor.Collapse(0)\n$insertStart = $anchor.Start\n\n$sentence = \" As\u00ed mismo, la contrase\u00f1a del usuario ser\u00e1 encriptada utilizando la librer\u00eda BCrypt que proporciona Java, utilizando una \" + [char]0x2018 + \"sal\" + [char]0x2019 + \" para reforzar a\u00fan m\u00e1s la contrase\u00f1a a la hora de almacenarla en la base de datos.\"\n\n# Insert the plain text first ...\n$anchor.InsertAfter($sentence)\n$newRange = $d.Range($insertStart, $anchor.End)\n\n# ... then replace that exact span via InsertXML so the new run gets the\n# same sz/szCs (12pt) run formatting used throughout the paragraph - plain\n# Range.InsertAfter()/Font.Size do not let us set the complex-script size\n# (w:szCs) that the rest of the document's runs carry.\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body><w:p><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">$sentence</w:t></w:r></w:p></w:body>\n</w:document>\n</pkg:xmlData></pkg:part></pkg:package>\n\"@\n$newRange.InsertXML($ooxml)\n"}
